$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct floating point rounding drift on a handful of pre-existing values ---
$ws.Range("L57").Value = 1455967.872
$ws.Range("H59").Value = 1180532.864
$ws.Range("L59").Value = 1425340.16
$ws.Range("P60").Value = -1491473.024
$ws.Range("P61").Value = 515001.024
$ws.Range("P62").Value = -107739.992
$ws.Range("T62").Value = -117905
$ws.Range("H63").Value = -161770.016
$ws.Range("T63").Value = -409462.944
$ws.Range("T80").Value = 29815.992

# --- Clear stale placeholder cells (leftover zeros / concatenation artifacts) ---
# These columns belong to periods that this statement does not actually report,
# so the cells should be empty rather than holding a numeric 0 / computed value.
$ws.Range("P57:W57").ClearContents()
$ws.Range("P58:W58").ClearContents()
$ws.Range("W59").ClearContents()
$ws.Range("W60").ClearContents()
$ws.Range("W61").ClearContents()
$ws.Range("W62").ClearContents()
$ws.Range("W63").ClearContents()
$ws.Range("B64:P64").ClearContents()
$ws.Range("W64").ClearContents()
$ws.Range("W65").ClearContents()
$ws.Range("W66").ClearContents()
$ws.Range("W67").ClearContents()
$ws.Range("W68").ClearContents()
$ws.Range("W69").ClearContents()
$ws.Range("W70").ClearContents()
$ws.Range("P71:W71").ClearContents()
$ws.Range("P72:W72").ClearContents()
$ws.Range("P73:W73").ClearContents()
$ws.Range("W74").ClearContents()
$ws.Range("W75").ClearContents()
$ws.Range("W76").ClearContents()
$ws.Range("P77:W77").ClearContents()
$ws.Range("P78:W78").ClearContents()
$ws.Range("W79").ClearContents()
$ws.Range("W80").ClearContents()
